$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1674.9
$ws.Range("I19").Value = 1288.3334
$ws.Range("J19").Value = 1840.5714
$ws.Range("K19").Value = 1288.3334
$ws.Range("L19").Value = 1840.5714
$ws.Range("M19").Value = -1113.3334
$ws.Range("N19").Value = -2190.5714
$ws.Range("H33").Value = 181.23529
$ws.Range("I33").Value = 190.0625
$ws.Range("K33").Value = 190.0625
$ws.Range("M33").Value = 38.9375
$ws.Range("H51").Value = 6857
$ws.Range("I51").Value = 5666
$ws.Range("K51").Value = 5666
$ws.Range("M51").Value = -5182
$ws.Range("H64").Value = 7060.091
$ws.Range("I64").Value = 5554
$ws.Range("J64").Value = 7210.7
$ws.Range("K64").Value = 5554
$ws.Range("L64").Value = 7210.7
$ws.Range("M64").Value = -5306
$ws.Range("N64").Value = -7706.7
$ws.Range("H67").Value = 7060.091
$ws.Range("I67").Value = 5554
$ws.Range("J67").Value = 7210.7
$ws.Range("K67").Value = 5554
$ws.Range("L67").Value = 7210.7
$ws.Range("M67").Value = -4696
$ws.Range("N67").Value = -8926.700000000001
$ws.Range("H70").Value = 92915.82000000001
$ws.Range("J70").Value = 144842
$ws.Range("L70").Value = 434526
$ws.Range("N70").Value = -435066
$ws.Range("H73").Value = 92915.82000000001
$ws.Range("J73").Value = 144842
$ws.Range("L73").Value = 434526
$ws.Range("N73").Value = -436398
$ws.Range("H76").Value = 66744572
$ws.Range("I76").Value = 98417.09
$ws.Range("K76").Value = 98417.09
$ws.Range("M76").Value = -98102.09
$ws.Range("H79").Value = 66744572
$ws.Range("I79").Value = 98417.09
$ws.Range("K79").Value = 98417.09
$ws.Range("M79").Value = -97325.09
$ws.Range("H107").Value = 39131.92
$ws.Range("I107").Value = 48295.715
$ws.Range("K107").Value = 48295.715
$ws.Range("M107").Value = -46375.715
$ws.Range("H129").Value = 13653.211
$ws.Range("I129").Value = 10161.538
$ws.Range("J129").Value = 21218.5
$ws.Range("K129").Value = 30484.614
$ws.Range("L129").Value = 63655.5
$ws.Range("M129").Value = -25484.614
$ws.Range("N129").Value = -73655.5
$ws.Range("H132").Value = 2392.9614
$ws.Range("I132").Value = 1392.7142
$ws.Range("K132").Value = 4178.142599999999
$ws.Range("M132").Value = -1648.142599999999
$ws.Range("H138").Value = 4110.27
$ws.Range("I138").Value = 2367.7778
$ws.Range("J138").Value = 5417.1387
$ws.Range("K138").Value = 7103.3334
$ws.Range("L138").Value = 16251.4161
$ws.Range("M138").Value = -1963.3334
$ws.Range("N138").Value = -26531.4161

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 272452.75
$ws.Range("I61").Value = 1291.3334
$ws.Range("K61").Value = 1291.3334
$ws.Range("M61").Value = -1079.3334
$ws.Range("H112").Value = 40257
$ws.Range("J112").Value = 40257
$ws.Range("L112").Value = 40257
$ws.Range("N112").Value = -43211
$ws.Range("H132").Value = 6975.5684
$ws.Range("I132").Value = 4497.5
$ws.Range("K132").Value = 13492.5
$ws.Range("M132").Value = -10962.5
$ws.Range("H136").Value = 272452.75
$ws.Range("I136").Value = 1291.3334
$ws.Range("K136").Value = 3874.0002
$ws.Range("M136").Value = -1324.0002

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("M75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("M78").ClearContents()
$ws.Range("H96").Value = 34499.5
$ws.Range("I96").Value = 34499.5
$ws.Range("K96").Value = 34499.5
$ws.Range("M96").Value = -31753.5

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5125.5
$ws.Range("I31").Value = 3945.889
$ws.Range("J31").Value = 6642.143
$ws.Range("K31").Value = 3945.889
$ws.Range("L31").Value = 6642.143
$ws.Range("M31").Value = -3650.889
$ws.Range("N31").Value = -7232.143
$ws.Range("H34").Value = 5125.5
$ws.Range("I34").Value = 3945.889
$ws.Range("J34").Value = 6642.143
$ws.Range("K34").Value = 3945.889
$ws.Range("L34").Value = 6642.143
$ws.Range("M34").Value = -3743.889
$ws.Range("N34").Value = -7046.143
$ws.Range("H94").Value = 4497
$ws.Range("J94").Value = 4497
$ws.Range("L94").Value = 4497
$ws.Range("N94").Value = -5399
$ws.Range("H132").Value = 5846.3823
$ws.Range("I132").Value = 5535.222
$ws.Range("K132").Value = 16605.666
$ws.Range("M132").Value = -14075.666
$ws.Range("H134").Value = 3591.2888
$ws.Range("I134").Value = 2910.121
$ws.Range("K134").Value = 8730.363000000001
$ws.Range("M134").Value = -6195.363000000001

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 3459.8
$ws.Range("J13").Value = 3974.75
$ws.Range("L13").Value = 11924.25
$ws.Range("N13").Value = -12260.25
$ws.Range("H139").Value = 944
$ws.Range("I139").Value = 938.4
$ws.Range("J139").Value = 1000
$ws.Range("K139").Value = 2815.2
$ws.Range("L139").Value = 3000
$ws.Range("M139").Value = 2324.8
$ws.Range("N139").Value = -13280
$ws.Range("H140").Value = 1099.4736
$ws.Range("I140").Value = 743.94446
$ws.Range("K140").Value = 2231.83338
$ws.Range("M140").Value = 2948.16662

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 168.07143
$ws.Range("I2").Value = 173.3077
$ws.Range("K2").Value = 173.3077
$ws.Range("M2").Value = -60.30770000000001
$ws.Range("H122").Value = 8118.3477
$ws.Range("I122").Value = 7151.4375
$ws.Range("J122").Value = 10328.429
$ws.Range("K122").Value = 21454.3125
$ws.Range("L122").Value = 30985.287
$ws.Range("M122").Value = -19004.3125
$ws.Range("N122").Value = -35885.287

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H6").Value = 215240
$ws.Range("J6").Value = 215240
$ws.Range("L6").Value = 215240
$ws.Range("N6").Value = -215464
$ws.Range("H92").Value = 45000
$ws.Range("J92").Value = 45000
$ws.Range("L92").Value = 45000
$ws.Range("N92").Value = -49992

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 45459520
$ws.Range("I122").Value = 76926580
$ws.Range("J122").Value = 7110.5557
$ws.Range("K122").Value = 230779740
$ws.Range("L122").Value = 21331.6671
$ws.Range("M122").Value = -230777290
$ws.Range("N122").Value = -26231.6671
$ws.Range("I136").Value = 3335166.2
$ws.Range("K136").Value = 10005498.6
$ws.Range("M136").Value = -10002948.6
